$wb = $excel.ActiveWorkbook

# The "DataTypes" worksheet holds the type-mapping table (Type, C#, MySQL, SQLite, T-SQL).
# A MySQL column (C) template is being started: fill in MySQL type names for each row,
# mirroring the existing SQLite (D) column values but using MySQL-specific spellings.
$ws = $wb.Worksheets.Item("DataTypes")

$ws.Range("C2").Value  = "BOOL"
$ws.Range("C3").Value  = "BLOB"
$ws.Range("C4").Value  = "CHAR"
$ws.Range("C5").Value  = "TEXT"
$ws.Range("C6").Value  = "VARCHAR(0)"
$ws.Range("C7").Value  = "INT"
$ws.Range("C8").Value  = "BIGINT"
$ws.Range("C9").Value  = "MEDIUMINT"
$ws.Range("C10").Value = "SMALLINT"
$ws.Range("C11").Value = "FLOAT"
$ws.Range("C12").Value = "DOUBLE"
$ws.Range("C13").Value = "DECIMAL(0, 0)"
$ws.Range("C14").Value = "DATE"
$ws.Range("C15").Value = "TIME"
$ws.Range("C16").Value = "DATETIME"

# Move the active selection, matching where the editor ended up on this sheet.
$ws.Range("H6").Select()
